$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title date line
Replace-Text "2024-04-09 Tuesday" "2024-04-10 Wednesday"

# Table of division problems (5x5), using table cell access to avoid
# any cross-replacement collisions between old/new values.
$tbl = $d.Tables.Item(1)

$values = @(
    @("63÷9=","50÷4="), @("49÷6=","57÷8="), @("46÷2=","78÷6="), @("71÷5=","19÷3="), @("19÷4=","34÷6="),
    @("10÷5=","15÷5="), @("70÷6=","79÷4="), @("78÷4=","64÷6="), @("92÷6=","12÷4="), @("96÷6=","14÷3="),
    @("32÷7=","21÷6="), @("26÷2=","95÷5="), @("54÷7=","24÷4="), @("40÷8=","37÷5="), @("48÷2=","51÷8="),
    @("75÷2=","84÷6="), @("57÷7=","38÷3="), @("15÷4=","70÷4="), @("40÷7=","12÷3="), @("32÷6=","34÷6="),
    @("60÷2=","56÷3="), @("90÷2=","49÷8="), @("34÷9=","35÷4="), @("61÷3=","36÷2="), @("37÷5=","97÷3=")
)

$contentRows = @(1, 5, 9, 13, 17)
$idx = 0
foreach ($r in $contentRows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($r, $c)
        $pair = $values[$idx]
        $old = $pair[0]
        $new = $pair[1]
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $new
        $idx = $idx + 1
    }
}
